# chore: adapt column header formatting to respective input file names (#7)
#
# Renames the "_old"/"_new" header-name suffixes to the respective
# format-version suffixes ("_FV2404" / "_FV2410"), wraps the sheet's used
# range in an Excel Table (ListObject) whose column names mirror the new
# headers, and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (row 1) -------------------------------
# Columns A:J -> "..._old"  => "..._FV2404"
# Column  K   -> "diff"     => unchanged
# Columns L:U -> "..._new"  => "..._FV2410"

$ws.Range("A1").Value = "Segmentname_FV2404"
$ws.Range("B1").Value = "Segmentgruppe_FV2404"
$ws.Range("C1").Value = "Segment_FV2404"
$ws.Range("D1").Value = "Datenelement_FV2404"
$ws.Range("E1").Value = "Segment ID_FV2404"
$ws.Range("F1").Value = "Code_FV2404"
$ws.Range("G1").Value = "Qualifier_FV2404"
$ws.Range("H1").Value = "Beschreibung_FV2404"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("J1").Value = "Bedingung_FV2404"

$ws.Range("L1").Value = "Segmentname_FV2410"
$ws.Range("M1").Value = "Segmentgruppe_FV2410"
$ws.Range("N1").Value = "Segment_FV2410"
$ws.Range("O1").Value = "Datenelement_FV2410"
$ws.Range("P1").Value = "Segment ID_FV2410"
$ws.Range("Q1").Value = "Code_FV2410"
$ws.Range("R1").Value = "Qualifier_FV2410"
$ws.Range("S1").Value = "Beschreibung_FV2410"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("U1").Value = "Bedingung_FV2410"

# --- 2. Wrap A1:U70 in an Excel Table (ListObject) ---------------------
$dataRange = $ws.Range("A1:U70")
$table = $ws.ListObjects.Add(1, $dataRange, $false, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row -------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "done"
